$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Move Robot2 to location (2, 8) and remove the toolkit."
$ws.Range("B2").Value = 36.320925
$ws.Range("C2").Value = 3971
$ws.Range("D2").Value = "'0.00852"
$ws.Range("E2").Value = "760e41a3-d8c2-444d-9572-7366ad81dec9"

# Row 3
$ws.Range("A3").Value = "Move Robot26 to location (11, 4) and remove the liquid spill."
$ws.Range("B3").Value = 46.797177
$ws.Range("C3").Value = 4009
$ws.Range("D3").Value = "'0.00906"
$ws.Range("E3").Value = "3d2c6b48-0f96-4493-8462-65b457c992b8"

# Row 4
$ws.Range("A4").Value = "Move Robot42 to location (9, 5) and remove the large debris."
$ws.Range("B4").Value = 61.136099
$ws.Range("C4").Value = 4615
$ws.Range("D4").Value = "'0.00993"
$ws.Range("E4").Value = "eef0bfa4-adcc-482e-8ccc-96ea2203be96"

# Row 5
$ws.Range("A5").Value = "Move Robot48 to location (5, 6) and remove the dust."
$ws.Range("B5").Value = 69.949413
$ws.Range("C5").Value = 4567
$ws.Range("D5").Value = "'0.00996"
$ws.Range("E5").Value = "b4bd5fd8-7da8-4975-a5ca-48fdff2fe5bf"

# Row 6
$ws.Range("A6").Value = "Move Robot31 to location (9, 4) and remove the grass."
$ws.Range("B6").Value = 35.20186
$ws.Range("C6").Value = 3866
$ws.Range("D6").Value = "'0.00828"
$ws.Range("E6").Value = "052d2de2-2857-4277-9766-9e3db72e5584"

# Row 7
$ws.Range("A7").Value = "Move Robot8 to location (8, 12) and remove the small debris."
$ws.Range("B7").Value = 42.525113
$ws.Range("C7").Value = 4484
$ws.Range("D7").Value = "'0.00915"
$ws.Range("E7").Value = "4f611a9a-7cef-4415-9d8f-88f33aded40b"

# Row 8
$ws.Range("A8").Value = "Move Robot23 to location (11, 1) and remove the vehicle."
$ws.Range("B8").Value = 32.039849
$ws.Range("C8").Value = 3865
$ws.Range("D8").Value = "'0.0081"
$ws.Range("E8").Value = "0829af54-b795-4c0b-b3d9-06230bbe514d"

# Row 9
$ws.Range("A9").Value = "Move Robot23 to location (12, 10) and remove the construction materials."
$ws.Range("B9").Value = 27.348842
$ws.Range("C9").Value = 3927
$ws.Range("D9").Value = "'0.00807"
$ws.Range("E9").Value = "eebbb157-65d7-4371-915d-d0e144c29e9f"

# Row 10
$ws.Range("A10").Value = "Move Robot14 to location (7, 11) and remove the tree branches."
$ws.Range("B10").Value = 26.50084
$ws.Range("C10").Value = 3848
$ws.Range("E10").Value = "c0d299d7-8690-4a94-9213-8e6aa34cce3e"

# Row 11
$ws.Range("A11").Value = "Move Robot15 to location (5, 3) and remove the screws."
$ws.Range("B11").Value = 26.786554
$ws.Range("C11").Value = 3843
$ws.Range("D11").Value = "'0.00813"
$ws.Range("E11").Value = "d03d844b-dbc4-4b45-9272-83187249b2a7"
